$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells: "<Name>_old" -> "<Name>_FV2310", "<Name>_new" -> "<Name>_FV2404"
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value2
    if ($val -ne $null) {
        if ($val.EndsWith("_old")) {
            $cell.Value2 = $val.Substring(0, $val.Length - 4) + "_FV2310"
        } elseif ($val.EndsWith("_new")) {
            $cell.Value2 = $val.Substring(0, $val.Length - 4) + "_FV2404"
        }
    }
}

# Turn the data range into an Excel Table ("Table1") with a header row
$range = $ws.Range("A1:U77")
$listObject = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $range, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$listObject.Name = "Table1"
$listObject.TableStyle = ""

# Freeze the header row
$ws.Application.ActiveWindow.SplitRow = 1
$ws.Application.ActiveWindow.FreezePanes = $true
